$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update adds a new "Primera"/"Segunda" price pair for this
# market (date 44491) at the top of the recent-history block, pushing the
# previously-newest rows (old 71-74) down by two rows (to 73-76).
$ws.Rows.Item(71).Insert()
$ws.Rows.Item(71).Insert()

# New row 71: Primera, week of 44491
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(71, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value = 44491
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 100112042
$ws.Cells.Item(71, 7).Value = "Locoto"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 160
$ws.Cells.Item(71, 11).Value = 29000
$ws.Cells.Item(71, 12).Value = 30000
$ws.Cells.Item(71, 13).Value = 29500
$ws.Cells.Item(71, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 1475
$ws.Cells.Item(71, 17).Value = 20
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# New row 72: Segunda, week of 44491
$ws.Cells.Item(72, 1).Value = 1
$ws.Cells.Item(72, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value = 44491
$ws.Cells.Item(72, 5).Value = 15
$ws.Cells.Item(72, 6).Value = 100112042
$ws.Cells.Item(72, 7).Value = "Locoto"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Segunda"
$ws.Cells.Item(72, 10).Value = 160
$ws.Cells.Item(72, 11).Value = 26000
$ws.Cells.Item(72, 12).Value = 27000
$ws.Cells.Item(72, 13).Value = 26500
$ws.Cells.Item(72, 14).Value = "`$/caja 20 kilos"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 1325
$ws.Cells.Item(72, 17).Value = 20
$ws.Cells.Item(72, 18).Value = "Hortaliza"
